# "added output for the engineer"
# Updates the first/last generation individual values (columns B and D)
# and the derived fitness/penalty statistics (G1, G2, B21, D21) with a
# freshly produced run's output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B - "first generation:" individuals
$ws.Range("B1").Value  = 7
$ws.Range("B2").Value  = 0
$ws.Range("B3").Value  = 12
$ws.Range("B4").Value  = 14
$ws.Range("B5").Value  = 15
$ws.Range("B6").Value  = 19
$ws.Range("B7").Value  = 9
$ws.Range("B8").Value  = 17
$ws.Range("B9").Value  = 18
$ws.Range("B10").Value = 5
$ws.Range("B11").Value = 8
$ws.Range("B12").Value = 11
$ws.Range("B13").Value = 4
$ws.Range("B14").Value = 10
$ws.Range("B15").Value = 3
$ws.Range("B16").Value = 16
$ws.Range("B17").Value = 1
$ws.Range("B18").Value = 13
$ws.Range("B19").Value = 2
$ws.Range("B20").Value = 6

# Column D - "last generation:" individuals
$ws.Range("D1").Value  = 1
$ws.Range("D2").Value  = 0
$ws.Range("D3").Value  = 7
$ws.Range("D4").Value  = 12
$ws.Range("D5").Value  = 3
$ws.Range("D6").Value  = 14
$ws.Range("D7").Value  = 4
$ws.Range("D8").Value  = 19
$ws.Range("D9").Value  = 5
$ws.Range("D10").Value = 15
$ws.Range("D11").Value = 6
$ws.Range("D12").Value = 9
$ws.Range("D13").Value = 13
$ws.Range("D14").Value = 16
$ws.Range("D15").Value = 2
$ws.Range("D16").Value = 11
$ws.Range("D17").Value = 10
$ws.Range("D18").Value = 8
$ws.Range("D19").Value = 18
$ws.Range("D20").Value = 17

# Fitness / penalty improvement summary values
$ws.Range("G1").Value  = 114.1620489728746
$ws.Range("G2").Value  = 67.63156526669685
$ws.Range("B21").Value = 0.6956393348003318
$ws.Range("D21").Value = 0.7941561180693337
